$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''245.47'
$ws.Range("E2").Value = '''-0.59%'

# Row 3
$ws.Range("D3").Value = '''27.16'
$ws.Range("E3").Value = '''3.32%'

# Row 4
$ws.Range("D4").Value = '''5.111'
$ws.Range("E4").Value = '''0.70%'

# Row 5
$ws.Range("D5").Value = '''0.05699'
$ws.Range("E5").Value = '''1.66%'

# Row 6
$ws.Range("D6").Value = '''6.511'
$ws.Range("E6").Value = '''0.44%'

# Row 7
$ws.Range("D7").Value = '''0.8192'
$ws.Range("E7").Value = '''0.75%'

# Row 8
$ws.Range("D8").Value = '''0.8596'
$ws.Range("E8").Value = '''1.98%'

# Row 9
$ws.Range("B9").Value = 'MandalaExchangeToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D9").Value = '''0.06947'
$ws.Range("E9").Value = '''-0.88%'

# Row 10
$ws.Range("B10").Value = 'BitrueCoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D10").Value = '''0.02852'
$ws.Range("E10").Value = '''0.24%'

# Row 11
$ws.Range("B11").Value = 'BitMartToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D11").Value = '''0.09394'
$ws.Range("E11").Value = '''-0.25%'

# Row 12
$ws.Range("B12").Value = 'BitForexToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D12").Value = '''0.001524'
$ws.Range("E12").Value = '''0.43%'

# Row 13
$ws.Range("B13").Value = 'CoinExToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D13").Value = '''0.04046'
$ws.Range("E13").Value = '''-13.04%'

# Row 14
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").Value = '''0.0006007'
$ws.Range("E14").Value = '''-0.05%'

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006212'
$ws.Range("E15").Value = '''0.63%'

# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.512'
$ws.Range("E16").Value = '''-2.65%'

# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''3.009'
$ws.Range("E17").Value = '''-0.16%'

# Row 18
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.230'
$ws.Range("E18").Value = '''8.48%'

# Row 19
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '''0.3165'
$ws.Range("E19").Value = '''1.23%'

# Row 20
$ws.Range("B20").Value = 'WazirX'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D20").Value = '''0.1332'
$ws.Range("E20").Value = '''-0.31%'

# Row 21
$ws.Range("D21").Value = '''0.03222'
$ws.Range("E21").Value = '''0.91%'

# Row 22
$ws.Range("E22").Value = '''-1.80%'

# Row 23
$ws.Range("D23").Value = '''3.572'
$ws.Range("E23").Value = '''-4.75%'

# Row 24
$ws.Range("E24").Value = '''1.74%'

# Row 25
$ws.Range("D25").Value = '''0.001218'
$ws.Range("E25").Value = '''-2.37%'

# Row 26
$ws.Range("E26").Value = '''-2.42%'

# Row 27
$ws.Range("D27").Value = '''0.00009897'
$ws.Range("E27").Value = '''3.12%'

# Row 28
$ws.Range("E28").Value = '''-25.27%'

# Row 40
$ws.Range("D40").Value = '''0.03731'
$ws.Range("E40").Value = '''1.76%'

# Row 41
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '''0.1059'
$ws.Range("E41").Value = '''-21.61%'

# Row 42
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '''0.002459'
$ws.Range("E42").Value = '''-7.53%'

# Row 43
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = '''0.003447'
$ws.Range("E43").Value = '''-43.94%'

# Row 44
$ws.Range("D44").Value = '''0.009712'
$ws.Range("E44").Value = '''17.62%'

# Row 45
$ws.Range("D45").Value = '''0.00005135'
$ws.Range("E45").Value = '''-4.46%'

# Row 46
$ws.Range("D46").Value = '''0.00000000750'
$ws.Range("E46").Value = '''-0.05%'

# Row 47
$ws.Range("D47").Value = '''0.1010'
$ws.Range("E47").Value = '''-8.23%'

# Row 48
$ws.Range("D48").Value = '''0.002505'
$ws.Range("E48").Value = '''-3.46%'

# Row 49
$ws.Range("D49").Value = '''0.00002099'
$ws.Range("E49").Value = '''-0.05%'

# Row 50
$ws.Range("D50").Value = '''0.0001999'
$ws.Range("E50").Value = '''-0.05%'
